# Terms of use: drop limitation on 500 installations per year
# https://phabricator.endlessm.com/T33418

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Body text: remove the "500 computers per year" installation cap and
#    replace it with "one or more devices".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "on 500 computers per year for", $true, $false, $false, $false, $false,
    $true, 1, $false, "on one or more devices for", 2)

# ---------------------------------------------------------------------------
# 2. Header: bump the "Last Updated" date from 26 May 2021 to 25 April 2022.
#    Apply the edits run-by-run (right-to-left so earlier offsets stay
#    valid) so the existing run/formatting structure is preserved as
#    closely as possible.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Sections.Count; $i++) {
    $sec = $d.Sections.Item($i)
    $hdr = $sec.Headers.Item(1)
    $full = $hdr.Range

    if ($full.Text.StartsWith("Last Updated: 26 May 2021")) {
        # "2021" -> "2022"
        $year = $full.Duplicate
        $year.SetRange(21, 25)
        $year.Text = "2022"

        # "May " -> "April "
        $month = $full.Duplicate
        $month.SetRange(17, 21)
        $month.Text = "April "

        # "6" -> "5"  (together with the preceding run's trailing "2" this
        # turns "Last Updated: 26" into "Last Updated: 25")
        $day = $full.Duplicate
        $day.SetRange(15, 16)
        $day.Text = "5"
    }
}
